# Update the "dayrealmob" forecast sheet: new date header and refreshed hourly values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "21/03/2023"

$values = @{
    2  = 89
    3  = 180
    4  = 235
    5  = 229
    6  = 193
    7  = 159
    8  = 180
    9  = 156
    10 = 165
    11 = 175
    12 = 167
    13 = 140
    14 = 83
    15 = 52
    16 = 25
    17 = 18
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
